# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit diff:
#  - Timestamp footer text bumped from 12:22 to 12:52
#  - Uzbekistan overtakes Eslovaquia/Crucero/Republica de Macedonia (rows 77-80)
#  - Etiopia overtakes Barbados (rows 139-140)
#  - Misc numeric refreshes for Pakistan (row 36), Kazajistan (row 75) and Malta (row 102)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 12:52"

# --- Row 36: Pakistan ---
$ws.Range("B36").Value = 4892
$ws.Range("C36").Value = 197
$ws.Range("E36").Value = 4053
$ws.Range("G36").Value = 11
$ws.Range("H36").Value = 77

# --- Row 75: Kazajistan ---
$ws.Range("D75").Value = 81
$ws.Range("E75").Value = 768

# --- Row 77: now Uzbekistan (was Eslovaquia) ---
$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("B77").Value = 729
$ws.Range("C77").Value = 105
$ws.Range("D77").Value = 42
$ws.Range("E77").Value = 684
$ws.Range("F77").Value = 8
$ws.Range("H77").Value = 3

# --- Row 78: now Eslovaquia (was Crucero) ---
$ws.Range("A78").Value = "Eslovaquia"
$ws.Range("B78").Value = 728
$ws.Range("C78").Value = 13
$ws.Range("D78").Value = 23
$ws.Range("E78").Value = 703
$ws.Range("F78").Value = 5
$ws.Range("H78").Value = 2

# --- Row 79: now Crucero (was Republica de Macedonia) ---
$ws.Range("A79").Value = "Crucero"
$ws.Range("B79").Value = 712
$ws.Range("D79").Value = 619
$ws.Range("E79").Value = 82
$ws.Range("F79").Value = 10
$ws.Range("H79").Value = 11

# --- Row 80: now Republica de Macedonia (was Uzbekistan) ---
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 711
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 41
$ws.Range("E80").Value = 638
$ws.Range("F80").Value = 15
$ws.Range("H80").Value = 32

# --- Row 102: Malta ---
$ws.Range("B102").Value = 370
$ws.Range("C102").Value = 20
$ws.Range("E102").Value = 351

# --- Row 139: now Etiopia (was Barbados) ---
$ws.Range("A139").Value = "Etiopia"
$ws.Range("B139").Value = 69
$ws.Range("C139").Value = 4
$ws.Range("D139").Value = 10
$ws.Range("E139").Value = 56
$ws.Range("F139").Value = 0
$ws.Range("H139").Value = 3

# --- Row 140: now Barbados (was Etiopia) ---
$ws.Range("A140").Value = "Barbados"
$ws.Range("B140").Value = 67
$ws.Range("D140").Value = 11
$ws.Range("E140").Value = 52
$ws.Range("F140").Value = 4
$ws.Range("H140").Value = 4
